$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'pantalones ala moda',
    'pantalones de basketball',
    'pantalones de compresion para mujer',
    'pantalones de licra',
    'pantalones de los 70',
    'panties 360',
    'panties barely there',
    'panties for women barely there',
    'panties support',
    'panties to hold in tummy',
    'pants after abdominal surgery',
    'pants capri for women',
    'pants exercise women',
    'pants for cold weather',
    'pants for gym',
    'pants for pilates',
    'pants for running',
    'pants for women black',
    'pants for women cold',
    'pants for women high waist',
    'pants for women tall',
    'pants high waist',
    'pants leggings for women',
    'pants soccer',
    'pants sport women',
    'pants three quarter',
    'pants tights for women',
    'pants warm women',
    'pants women',
    'pants women capri',
    'pants women training',
    'pants womens',
    'pants xxs',
    'pantsfor women',
    'panty compression',
    'panty hose game',
    'panty hose knee highs',
    'panty hose leggings',
    'panty hose seam',
    'panty hose with back seam',
    'panty hose with seams',
    'panty tights',
    'pantyhose compression women',
    'pantyhose material',
    'pantyhose shorts',
    'pantyhose warm skin',
    'pantys vs',
    'paragon fit wear',
    'paragon knee',
    'paragon shoes',
    'paragon shoes for man',
    'paragon sports',
    'paragon workout clothes',
    'pastel athletic leggings',
    'pastel blue leggings',
    'pastel clothing',
    'pastel dress pants men',
    'pastel felt fabric',
    'pastel green leggings',
    'pastel leg warmers',
    'pastel leggings',
    'pastel mens pants',
    'pastel pants women',
    'pastel rainbow',
    'pastel rainbow leggings',
    'pastel rainbow leggings women',
    'pastel rainbow patch',
    'pastel tape',
    'pastel tight',
    'pastel tights',
    'pastel workout leggings',
    'patch for pants hole',
    'patched pants',
    'patella brace for women',
    'patella compression',
    'patella tracking',
    'patella tracking knee brace',
    'patella tracking knee brace for women',
    'pelvic brace for men',
    'pelvic brace for women',
    'pelvic compression shorts',
    'pelvic pain supplement',
    'pelvic support brace',
    'pelvic support for women',
    'pelvis brace',
    'pelvis support',
    'pelvis support for women',
    'pennies for basketball',
    'pennies for soccer',
    'penny ring',
    'performance compression tights',
    'performance gear',
    'performance grip basketball',
    'performance leggings women',
    'performance pants women',
    'performance waist pack',
    'period shorts',
    'period shorts for women',
    'personal generator',
    'personal inhalers'
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
